$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Re-colour the presentation's (shared) theme from "Integral" /
#    Red-Violet over to the stock "Office Theme" colour palette.
#    The font scheme and format scheme are already identical between
#    the two themes in this deck, so only the 12 theme colours need
#    to change. Writing through Slide.ThemeColorScheme updates the
#    theme part shared by the slide master (and therefore every
#    slide), exactly like choosing a different set of Theme Colors
#    from the Design tab in the real PowerPoint UI.
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72

# ------------------------------------------------------------------
# 2) Switch the table on slide 5 over to the built-in table style
#    that the author picked (GUID taken straight from the target
#    OOXML's <a:tableStyleId>).
# ------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{99A806F2-76C6-40FC-8565-F9F6E2745A49}")
